# Avance reporte de actividades y metas
# Adds a helper column (I) that builds a SQL "insert" statement for each
# data row, concatenating the CLASE (A), DEP (B), ACTIVIDAD (D), META (E)
# and CATEGORIA PRESUPUESTAL (G) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Last data row (sheet has header in row 1, data in rows 2:96)
$lastRow = 96

$formula = "=CONCATENATE(""insert import(clase, dep, actividad, meta, cate)values('"",A2,""', '"",B2,""', '"",D2,""', '"",E2,""', '"",G2,""')"")"
$ws.Range("I2:I$lastRow").Formula = $formula

# Re-touch the widths of columns C, F and H (visual layout tweak that came
# together with the new column being added).
$ws.Columns.Item(3).ColumnWidth = 52.666666666666664
$ws.Columns.Item(6).ColumnWidth = 22.833333333333332
$ws.Columns.Item(8).ColumnWidth = 46.666666666666664
